$wb = $excel.ActiveWorkbook

# --- Create the new "FSAE_Achilles" sheet by cloning "Trailer_Kumanzi" ---
$srcSheet = $wb.Worksheets.Item("Trailer_Kumanzi")
$srcSheet.Copy($null, $srcSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "FSAE_Achilles"

# --- Update the numeric / text content for the new template ---
$newSheet.Range("H3").Value = "FSAE_Achilles"

$newSheet.Range("F6").Value = -1.53

$newSheet.Range("F7").Value = -0.8
$newSheet.Range("G7").Value = 0
$newSheet.Range("H7").Value = 0.289

$newSheet.Range("F8").Value = -1

$newSheet.Range("F9").Value = 0.25
$newSheet.Range("H9").Value = 0.403

$newSheet.Range("F10").Value = -1.75
$newSheet.Range("H10").Value = 0.403

$newSheet.Range("H11").Formula = "=0.619*2+0.2"

$newSheet.Range("H12").Value = 165

$newSheet.Range("F13").Value = 43
$newSheet.Range("G13").Value = 192
$newSheet.Range("H13").Value = 206

# The "Trailer_Kumanzi" template carries a column of blank helper cells in
# column K (formatted but empty) plus a "Empty?" guess note; the new
# template only keeps the "guesses" note in K7, matching the other
# non-trailer sheets (e.g. Sedan_Hamba).
$newSheet.Range("K5:K13").Clear()
$newSheet.Range("K7").Value = "guesses"

# --- Fix up cell formatting so it matches the "highlighted input" styling
#     used elsewhere in the workbook (Sedan_Hamba uses the same visual style
#     for these kinds of input cells). ---
$fmtSrc = $wb.Worksheets.Item("Sedan_Hamba")

$fmtSrc.Range("F7:H7").Copy()
$newSheet.Range("F7:H7").PasteSpecial(-4122)

$fmtSrc.Range("F8").Copy()
$newSheet.Range("F8").PasteSpecial(-4122)

$fmtSrc.Range("K7").Copy()
$newSheet.Range("K7").PasteSpecial(-4122)

$fmtSrc.Range("H12").Copy()
$newSheet.Range("H12").PasteSpecial(-4122)

$fmtSrc.Range("F13:H13").Copy()
$newSheet.Range("F13:H13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Update the remembered selections on the two other sheets that changed
#     as part of the same editing session ---
$sedanHamba = $wb.Worksheets.Item("Sedan_Hamba")
$sedanHamba.Range("E22").Select()

$sedanHambaLG = $wb.Worksheets.Item("Sedan_HambaLG")
$sedanHambaLG.Range("H12").Select()

# --- Finally make the newly added sheet the active / selected tab, with its
#     own remembered selection ---
$newSheet.Range("G27").Select()
